$wb = $excel.ActiveWorkbook

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 174.94
$ws.Range("I15").Value = 174.94
$ws.Range("K15").Value = 524.8199999999999
$ws.Range("M15").Value = -355.8199999999999
$ws.Range("H86").Value = 4785.148
$ws.Range("I86").Value = 1029.1765
$ws.Range("J86").Value = 11170.3
$ws.Range("K86").Value = 1029.1765
$ws.Range("L86").Value = 11170.3
$ws.Range("M86").Value = 93.82349999999997
$ws.Range("N86").Value = -13416.3
$ws.Range("H88").Value = 492.05554
$ws.Range("I88").Value = 357.5
$ws.Range("J88").Value = 559.3333
$ws.Range("K88").Value = 357.5
$ws.Range("L88").Value = 559.3333
$ws.Range("M88").Value = 48.5
$ws.Range("N88").Value = -1371.3333
$ws.Range("H89").Value = 4785.148
$ws.Range("I89").Value = 1029.1765
$ws.Range("J89").Value = 11170.3
$ws.Range("K89").Value = 5145.8825
$ws.Range("L89").Value = 55851.5
$ws.Range("M89").Value = 470.1175000000003
$ws.Range("N89").Value = -67083.5
$ws.Range("H91").Value = 492.05554
$ws.Range("I91").Value = 357.5
$ws.Range("J91").Value = 559.3333
$ws.Range("K91").Value = 357.5
$ws.Range("L91").Value = 559.3333
$ws.Range("M91").Value = 1046.5
$ws.Range("N91").Value = -3367.3333
$ws.Range("H98").Value = 995.75
$ws.Range("I98").Value = 908
$ws.Range("J98").Value = 1142
$ws.Range("K98").Value = 908
$ws.Range("L98").Value = 1142
$ws.Range("M98").Value = 590
$ws.Range("N98").Value = -4138
$ws.Range("H100").Value = 2222.2222
$ws.Range("J100").Value = 2533.3333
$ws.Range("L100").Value = 2533.3333
$ws.Range("N100").Value = -3615.3333
$ws.Range("H122").Value = 995.75
$ws.Range("I122").Value = 908
$ws.Range("J122").Value = 1142
$ws.Range("K122").Value = 2724
$ws.Range("L122").Value = 3426
$ws.Range("M122").Value = -274
$ws.Range("N122").Value = -8326
$ws.Range("H129").Value = 257112.2
$ws.Range("J129").Value = 303808.94
$ws.Range("L129").Value = 911426.8200000001
$ws.Range("N129").Value = -921426.8200000001
$ws.Range("H137").Value = 99088.37
$ws.Range("I137").Value = 122362.73
$ws.Range("K137").Value = 367088.19
$ws.Range("M137").Value = -364538.19
$ws.Range("H139").Value = 50515
$ws.Range("J139").Value = 50515
$ws.Range("L139").Value = 50515
$ws.Range("N139").Value = -60795

# Sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 785.64703
$ws.Range("I2").Value = 704.3077
$ws.Range("K2").Value = 704.3077
$ws.Range("M2").Value = -591.3077
$ws.Range("H32").Value = 6707.8423
$ws.Range("I32").Value = 4944.681
$ws.Range("J32").Value = 24087.572
$ws.Range("K32").Value = 4944.681
$ws.Range("L32").Value = 24087.572
$ws.Range("M32").Value = -4657.681
$ws.Range("N32").Value = -24661.572
$ws.Range("H63").Value = 4466565.5
$ws.Range("I63").Value = 2692
$ws.Range("K63").Value = 2692
$ws.Range("M63").Value = -2006
$ws.Range("H66").Value = 4466565.5
$ws.Range("I66").Value = 2692
$ws.Range("K66").Value = 13460
$ws.Range("M66").Value = -10028
$ws.Range("H102").Value = 1569.7273
$ws.Range("I102").Value = 1474.1111
$ws.Range("K102").Value = 1474.1111
$ws.Range("M102").Value = 147.8888999999999
$ws.Range("H116").Value = 785.64703
$ws.Range("I116").Value = 704.3077
$ws.Range("K116").Value = 704.3077
$ws.Range("M116").Value = 1589.6923

# Sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 785.64703
$ws.Range("I3").Value = 704.3077
$ws.Range("K3").Value = 704.3077
$ws.Range("M3").Value = -590.3077
$ws.Range("H99").Value = 1554.1765
$ws.Range("I99").Value = 1540
$ws.Range("K99").Value = 1540
$ws.Range("M99").Value = -42
$ws.Range("H105").Value = 1251955.2
$ws.Range("I105").Value = 1671.1765
$ws.Range("J105").Value = 2176078.2
$ws.Range("K105").Value = 1671.1765
$ws.Range("L105").Value = 2176078.2
$ws.Range("M105").Value = 75.82349999999997
$ws.Range("N105").Value = -2179572.2
$ws.Range("H134").Value = 4024.853
$ws.Range("I134").Value = 4057.9033
$ws.Range("J134").Value = 3683.3333
$ws.Range("K134").Value = 12173.7099
$ws.Range("L134").Value = 11049.9999
$ws.Range("M134").Value = -9638.7099
$ws.Range("N134").Value = -16119.9999

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3749.9792
$ws.Range("I31").Value = 2438.9375
$ws.Range("J31").Value = 4405.5
$ws.Range("K31").Value = 2438.9375
$ws.Range("L31").Value = 4405.5
$ws.Range("M31").Value = -2143.9375
$ws.Range("N31").Value = -4995.5
$ws.Range("H34").Value = 3749.9792
$ws.Range("I34").Value = 2438.9375
$ws.Range("J34").Value = 4405.5
$ws.Range("K34").Value = 2438.9375
$ws.Range("L34").Value = 4405.5
$ws.Range("M34").Value = -2236.9375
$ws.Range("N34").Value = -4809.5
$ws.Range("H105").Value = 1086.4667
$ws.Range("I105").Value = 952.7692
$ws.Range("J105").Value = 1955.5
$ws.Range("K105").Value = 952.7692
$ws.Range("L105").Value = 1955.5
$ws.Range("M105").Value = 794.2308
$ws.Range("N105").Value = -5449.5
$ws.Range("H122").Value = 1009.7619
$ws.Range("I122").Value = 821.5454999999999
$ws.Range("K122").Value = 2464.6365
$ws.Range("M122").Value = -14.63649999999961

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 729.67
$ws.Range("I131").Value = 290.8
$ws.Range("J131").Value = 752.76843
$ws.Range("K131").Value = 872.4000000000001
$ws.Range("L131").Value = 2258.30529
$ws.Range("M131").Value = 4167.6
$ws.Range("N131").Value = -12338.30529

# Sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2551.3125
$ws.Range("I102").Value = 1942.3334
$ws.Range("J102").Value = 4378.25
$ws.Range("K102").Value = 1942.3334
$ws.Range("L102").Value = 4378.25
$ws.Range("M102").Value = -320.3334
$ws.Range("N102").Value = -7622.25
$ws.Range("H122").Value = 4499.75
$ws.Range("I122").Value = 3437.125
$ws.Range("K122").Value = 10311.375
$ws.Range("M122").Value = -7861.375
$ws.Range("H132").Value = 16957.135
$ws.Range("I132").Value = 3565.0322
$ws.Range("J132").Value = 86149.664
$ws.Range("K132").Value = 10695.0966
$ws.Range("L132").Value = 258448.992
$ws.Range("M132").Value = -8165.096600000001
$ws.Range("N132").Value = -263508.992

# Sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4560.375
$ws.Range("I7").Value = 4467.5884
$ws.Range("J7").Value = 4785.7144
$ws.Range("K7").Value = 4467.5884
$ws.Range("L7").Value = 4785.7144
$ws.Range("M7").Value = -4355.5884
$ws.Range("N7").Value = -5009.7144
$ws.Range("H68").Value = 3374.5
$ws.Range("J68").Value = 2999
$ws.Range("L68").Value = 2999
$ws.Range("N68").Value = -4497
$ws.Range("H71").Value = 3374.5
$ws.Range("J71").Value = 2999
$ws.Range("L71").Value = 14995
$ws.Range("N71").Value = -22483
$ws.Range("H100").Value = 1996.7333
$ws.Range("I100").Value = 994.6
$ws.Range("J100").Value = 2497.8
$ws.Range("K100").Value = 994.6
$ws.Range("L100").Value = 2497.8
$ws.Range("M100").Value = -453.6
$ws.Range("N100").Value = -3579.8
$ws.Range("H126").Value = 4560.375
$ws.Range("I126").Value = 4467.5884
$ws.Range("J126").Value = 4785.7144
$ws.Range("K126").Value = 13402.7652
$ws.Range("L126").Value = 14357.1432
$ws.Range("M126").Value = -10932.7652
$ws.Range("N126").Value = -19297.1432

# Sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1911.3438
$ws.Range("I126").Value = 1426.96
$ws.Range("J126").Value = 3641.2856
$ws.Range("K126").Value = 4280.88
$ws.Range("L126").Value = 10923.8568
$ws.Range("M126").Value = -1810.88
$ws.Range("N126").Value = -15863.8568
$ws.Range("H136").Value = 28676830
$ws.Range("I136").Value = 39703576
$ws.Range("J136").Value = 7289.5
$ws.Range("K136").Value = 119110728
$ws.Range("L136").Value = 21868.5
$ws.Range("M136").Value = -119108178
$ws.Range("N136").Value = -26968.5
